$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F (classes) ---
$ws.Range("F2").Value = 5
$ws.Range("F3").Value = 5
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = 4

# --- Column G (last_changes_of_class) ---
# New date-only number format is introduced: first the lowercase variant is
# tried (registers numFmt 166), then the uppercase variant is applied
# (registers numFmt 167 and is the one that actually sticks for the xf).
$ws.Range("G2:G5").NumberFormat = "yyyy-mm-dd"
$ws.Range("G2:G5").NumberFormat = "YYYY-MM-DD"
$ws.Range("G2").Value = 44966
$ws.Range("G3").Value = 44966
$ws.Range("G4").Value = 44966
$ws.Range("G5").Value = 44966

# --- Column H (date_becoming) ---
# Only row 3 changes value; keeps its existing number format/style.
$ws.Range("H3").Value = 44966

# --- Column I (recalling) ---
$ws.Range("I3").Value = 5
$ws.Range("I4").Value = 3
